# 7.2.1.1 Renewable energy share — add 2023 data column (T), update 2022
# "in per cent" -> "in percent" note, and correct the 2022 percentage figure.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the English note text (row 2, col C): "(in per cent)" -> "(in percent)"
$ws.Range("C2").Value = "(in percent)"

# --- Correct the 2022 renewable-share percentage (S5): 30 -> 29.9
$ws.Range("S5").Value = 29.9

# --- Add the new 2023 column (T): year header, percentage, hydro output
# Header (row 4) — copy S4's formatting onto T4, then set the year value.
$ws.Range("T4").Value = 2023
$ws.Range("S4").Copy()
$ws.Range("T4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Renewable-energy share row (row 5)
$ws.Range("T5").Value = 29.5
$ws.Range("S5").Copy()
$ws.Range("T5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# Hydropower production row (row 6)
$ws.Range("T6").Value = 12030.6
$ws.Range("S6").Copy()
$ws.Range("T6").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

# --- Match the new column width used across D:T (width 9)
$ws.Range("D1:T1").ColumnWidth = 8.2
